# #feat: escopo de materiais na familia de cavalete
#
# Insert a new line item into the "Religacao" family on the "unitario"
# sheet: "472000 / RESTABELECIDA LIG AGUA COM SERV ADIC", right above the
# existing "475500 / RESTABELECIDA LIG AGUA MUDAN TITULAR OP" row (row 46),
# pushing that row and everything below it down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("unitario")
$ws.Activate()

# Insert a blank row at 46 - this shifts the old row 46 (and everything
# after it) down to row 47, etc.
$ws.Rows.Item(46).Insert()

# The freshly inserted row loses some of the category-header formatting
# (border/fill) that Excel would normally carry over, so copy the format
# explicitly from the row that used to be the header (now at row 47,
# directly below the new row).
$ws.Range("A47:C47").Copy()
$ws.Range("A46:C46").PasteSpecial(-4122)
$ws.Rows.Item(46).RowHeight = 18.75

# Populate the new row with the new catalogue entry.
$ws.Cells.Item(46, 1).Value = "472000"
$ws.Cells.Item(46, 2).Value = "RESTABELECIDA LIG AGUA COM SERV ADIC"
$ws.Cells.Item(46, 3).Value = "Religacao"

# The _FilterDatabase range on "unitario" covers the header+data rows; it
# needs to grow by one row to keep including the (now shifted) last row.
$wb.Names.Item(1).RefersTo = "=unitario!`$A`$3:`$C`$56"

# Restore the on-screen scroll/selection state.
$excel.ActiveWindow.ScrollRow = 34
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C47").Select()
